$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "EmpID"
$ws.Range("A1:D1").Style = "Normal"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").RowHeight = 15
$ws.Range("A1:D1").Select() | Out-Null
